$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep numeric-looking Price values stored as text (column is text-typed)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated coin data
$ws.Range("D2").Value = '25.897.90'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '1.620.54'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").Value = '213.52'
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("D6").Value = '0.499'
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("D8").Value = '0.249'
$ws.Range("E8").Value = '  -1.01%  '

$ws.Range("D9").Value = '0.0616'
$ws.Range("E9").Value = '  -2.79%  '

$ws.Range("D10").Value = '18.29'
$ws.Range("E10").Value = '  -5.49%  '

$ws.Range("D11").Value = '0.0789'
$ws.Range("E11").Value = '  -0.40%  '

$ws.Range("D12").Value = '1.846.80'
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").Value = '1.624.55'
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D14").Value = '4.14'
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").Value = '0.523'
$ws.Range("E15").Value = '  -2.84%  '

$ws.Range("D16").Value = '25.903.18'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").Value = '0.0₃0739'
$ws.Range("E17").Value = '  -2.35%  '

$ws.Range("D18").Value = '61.35'
$ws.Range("E18").Value = '  -2.44%  '

$ws.Range("E19").Value = '  +0.43%  '

$ws.Range("D20").Value = '190.96'
$ws.Range("E20").Value = '  -0.77%  '

$ws.Range("D21").Value = '4.23'
$ws.Range("E21").Value = '  -1.14%  '

$ws.Range("D22").Value = '9.48'
$ws.Range("E22").Value = '  -2.52%  '

$ws.Range("D23").Value = '6.03'
$ws.Range("E23").Value = '  -1.48%  '

$ws.Range("D24").Value = '0.133'
$ws.Range("E24").Value = '  +1.73%  '

$ws.Range("D25").Value = '143.65'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("E26").Value = '  +0.63%  '

$ws.Range("D27").Value = '1.72'
$ws.Range("E27").Value = '  -2.67%  '

$ws.Range("D28").Value = '6.66'
$ws.Range("E28").Value = '  -2.15%  '

$ws.Range("D29").Value = '15.21'
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("E30").Value = '  -0.51%  '

$ws.Range("D31").Value = '0.0480'
$ws.Range("E31").Value = '  -1.71%  '

$ws.Range("D32").Value = '3.12'
$ws.Range("E32").Value = '  -3.40%  '

$ws.Range("D33").Value = '3.10'
$ws.Range("E33").Value = '  -4.67%  '

$ws.Range("E34").Value = '  -1.35%  '

$ws.Range("D35").Value = '1.48'
$ws.Range("E35").Value = '  -2.52%  '

$ws.Range("D36").Value = '1.119.56'
$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("D37").Value = '0.847'
$ws.Range("E37").Value = '  -5.24%  '

$ws.Range("D38").Value = '2.39'
$ws.Range("E38").Value = '  -2.79%  '

$ws.Range("D39").Value = '0.514'
$ws.Range("E39").Value = '  -3.41%  '

$ws.Range("E40").Value = '  -1.86%  '

$ws.Range("D41").Value = '97.90'
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("D42").Value = '1.756.79'
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").Value = '0.748'
$ws.Range("E43").Value = '  -5.58%  '

$ws.Range("D44").Value = '5.08'
$ws.Range("E44").Value = '  -4.40%  '

$ws.Range("D45").Value = '0.0₆0113'
$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '54.11'
$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.49'
$ws.Range("E47").Value = '  -0.53%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0521'
$ws.Range("E48").Value = '  +0.35%  '

$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '1.01'
$ws.Range("E50").Value = '  +0.63%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.46'
$ws.Range("E51").Value = '  -3.17%  '
